$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.NumberFormat = "General"
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '61.256.64'
Set-TextValue $ws.Range('E2') '  +0.86%  '
Set-TextValue $ws.Range('D3') '2.935.17'
Set-TextValue $ws.Range('E3') '  +1.01%  '
Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  +0.05%  '
Set-TextValue $ws.Range('D5') '592.86'
Set-TextValue $ws.Range('E5') '  +1.02%  '
Set-TextValue $ws.Range('D6') '145.97'
Set-TextValue $ws.Range('E6') '  +1.24%  '
Set-TextValue $ws.Range('D7') '1.00'
Set-TextValue $ws.Range('E7') '  +0.03%  '
Set-TextValue $ws.Range('E8') '  +0.77%  '
Set-TextValue $ws.Range('D9') '7.00'
Set-TextValue $ws.Range('E9') '  +4.77%  '
Set-TextValue $ws.Range('E10') '  +0.31%  '
Set-TextValue $ws.Range('E11') '  -0.33%  '
Set-TextValue $ws.Range('D12') '0.0000226'
Set-TextValue $ws.Range('E12') '  +0.59%  '
Set-TextValue $ws.Range('D13') '33.92'
Set-TextValue $ws.Range('E13') '  +1.60%  '
Set-TextValue $ws.Range('D15') '3.422.30'
Set-TextValue $ws.Range('E15') '  +1.09%  '
Set-TextValue $ws.Range('D16') '61.072.49'
Set-TextValue $ws.Range('E16') '  +0.58%  '
Set-TextValue $ws.Range('E17') '  +1.11%  '
Set-TextValue $ws.Range('D18') '2.937.39'
Set-TextValue $ws.Range('E18') '  +1.08%  '
Set-TextValue $ws.Range('D19') '437.65'
Set-TextValue $ws.Range('E19') '  +2.34%  '
Set-TextValue $ws.Range('E20') '  -0.20%  '
Set-TextValue $ws.Range('D21') '0.681'
Set-TextValue $ws.Range('E21') '  +0.09%  '
Set-TextValue $ws.Range('E22') '  +0.94%  '
Set-TextValue $ws.Range('D23') '81.72'
Set-TextValue $ws.Range('E23') '  +1.19%  '
Set-TextValue $ws.Range('D24') '11.00'
Set-TextValue $ws.Range('E24') '  +1.33%  '
Set-TextValue $ws.Range('E25') '  +0.01%  '
Set-TextValue $ws.Range('D26') '11.95'
Set-TextValue $ws.Range('E26') '  +0.89%  '
Set-TextValue $ws.Range('E27') '  +0.07%  '
Set-TextValue $ws.Range('D28') '2.27'
Set-TextValue $ws.Range('E28') '  +3.82%  '
Set-TextValue $ws.Range('E29') '  +0.42%  '
Set-TextValue $ws.Range('E30') '  -2.15%  '
Set-TextValue $ws.Range('E31') '  +3.88%  '
Set-TextValue $ws.Range('D32') '26.73'
Set-TextValue $ws.Range('E32') '  +1.40%  '
Set-TextValue $ws.Range('D33') '1.00'
Set-TextValue $ws.Range('E33') '  +0.05%  '
Set-TextValue $ws.Range('D34') '0.0₃0874'
Set-TextValue $ws.Range('E34') '  +2.04%  '
Set-TextValue $ws.Range('D35') '1.02'
Set-TextValue $ws.Range('E35') '  +0.74%  '
Set-TextValue $ws.Range('E36') '  +1.59%  '
Set-TextValue $ws.Range('E37') '  +0.31%  '
Set-TextValue $ws.Range('B38') 'Kaspa'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D38') '0.124'
Set-TextValue $ws.Range('E38') '  +0.91%  '
Set-TextValue $ws.Range('B39') 'Stacks'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D39') '2.00'
Set-TextValue $ws.Range('E39') '  +1.41%  '
Set-TextValue $ws.Range('D40') '8.62'
Set-TextValue $ws.Range('E40') '  +0.59%  '
Set-TextValue $ws.Range('D41') '42.25'
Set-TextValue $ws.Range('E41') '  +2.24%  '
Set-TextValue $ws.Range('D42') '0.288'
Set-TextValue $ws.Range('E42') '  -2.37%  '
Set-TextValue $ws.Range('D43') '377.67'
Set-TextValue $ws.Range('E43') '  +1.08%  '
Set-TextValue $ws.Range('E44') '  -0.48%  '
Set-TextValue $ws.Range('D45') '2.697.46'
Set-TextValue $ws.Range('E45') '  +0.29%  '
Set-TextValue $ws.Range('D46') '133.21'
Set-TextValue $ws.Range('E46') '  +0.86%  '
Set-TextValue $ws.Range('D48') '24.07'
Set-TextValue $ws.Range('E48') '  +0.36%  '
Set-TextValue $ws.Range('E49') '  -0.12%  '
Set-TextValue $ws.Range('E50') '  -1.03%  '
Set-TextValue $ws.Range('E51') '  +1.47%  '
